$d = $word.ActiveDocument
$d.Content.Find.Execute("explaratoire", $false, $false, $false, $false, $false, $true, 1, $false, "exploratoire", 2)
